$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Software Component")

# Add the new "Out of scope" column to the table (Table1), which also
# extends the table range (A1:C3 -> A1:D3) and the autofilter range.
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()

# Header + data for the new column.
$col.Range.Cells(1, 1).Value = "Out of scope"
$col.Range.Cells(2, 1).Value = "No"
$col.Range.Cells(3, 1).Value = "Yes"

# Match column width/bestfit behaviour of the other columns.
$ws.Columns.Item(4).AutoFit()
